# Datos Covid-19 Murcia 200614
# Fill in the two new daily rows (92, 93) that were previously blank
# IFERROR(...,"") placeholders, and refresh the view's frozen-pane
# scroll position / active selection to match where the editor was
# working (near the bottom of the data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 92 (14/06/2020) ---------------------------------------------------
$ws.Range("A92").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),43995)'
$ws.Range("B92").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),69)'
$ws.Range("C92").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),59)'
$ws.Range("D92").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),10)'
$ws.Range("E92").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),3)'
$ws.Range("F92").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),1413)'
$ws.Range("G92").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),150)'
$ws.Range("H92").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),3190)'
$ws.Range("A92").NumberFormat = "d/m/yyyy"
$ws.Range("I92").Formula = "=+H92-H91"
$ws.Range("J92").Formula = "=+F92-F91"

# --- Row 93 (15/06/2020) ---------------------------------------------------
$ws.Range("A93").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),43996)'
$ws.Range("B93").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),68)'
$ws.Range("C93").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),58)'
$ws.Range("D93").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),10)'
$ws.Range("E93").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),3)'
$ws.Range("F93").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),1414)'
$ws.Range("G93").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),150)'
$ws.Range("H93").Formula = '=IFERROR(__xludf.DUMMYFUNCTION("""COMPUTED_VALUE"""),3190)'
$ws.Range("A93").NumberFormat = "d/m/yyyy"
$ws.Range("I93").Formula = "=+H93-H92"
$ws.Range("J93").Formula = "=+F93-F92"

# --- View: frozen-pane scroll position + active selection ------------------
$ws.Range("F97").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 77
$win.ScrollColumn = 1
